$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "260.77"
Set-TextValue "E2" "1.73%"
Set-TextValue "D3" "27.12"
Set-TextValue "E3" "1.24%"
Set-TextValue "D4" "4.738"
Set-TextValue "E4" "5.82%"
Set-TextValue "D5" "0.06081"
Set-TextValue "E5" "3.42%"
Set-TextValue "D6" "6.665"
Set-TextValue "E6" "0.86%"
Set-TextValue "D7" "0.8474"
Set-TextValue "D8" "0.9216"
Set-TextValue "E8" "-0.75%"
Set-TextValue "D9" "0.1406"
Set-TextValue "E9" "2.13%"
Set-TextValue "D10" "0.04948"
Set-TextValue "E10" "8.55%"
Set-TextValue "E11" "0.82%"
Set-TextValue "D12" "0.03136"
Set-TextValue "E12" "2.30%"
Set-TextValue "D13" "0.09076"
Set-TextValue "E13" "-0.08%"
Set-TextValue "D14" "0.001546"
Set-TextValue "E14" "1.41%"
Set-TextValue "D15" "0.0006062"
Set-TextValue "E15" "0.37%"
Set-TextValue "D16" "0.006129"
Set-TextValue "E16" "-1.10%"
Set-TextValue "D17" "3.453"
Set-TextValue "E17" "-0.84%"
Set-TextValue "D18" "3.154"
Set-TextValue "E18" "-0.53%"
Set-TextValue "D19" "2.167"
Set-TextValue "E19" "-1.65%"
Set-TextValue "E20" "2.60%"
Set-TextValue "E21" "0.85%"
Set-TextValue "D22" "4.096"
Set-TextValue "E22" "4.45%"
Set-TextValue "D23" "0.04261"
Set-TextValue "E23" "0.15%"
Set-TextValue "D24" "0.001218"
Set-TextValue "E24" "-0.24%"
Set-TextValue "E25" "-9.11%"
Set-TextValue "E26" "0.05%"
Set-TextValue "D27" "0.0001574"
Set-TextValue "E27" "3.36%"
Set-TextValue "D40" "0.03876"
Set-TextValue "D41" "0.1113"
Set-TextValue "E41" "1.37%"
Set-TextValue "D42" "0.004129"
Set-TextValue "E42" "-34.34%"
Set-TextValue "E43" "18.11%"
Set-TextValue "E44" "0.39%"
Set-TextValue "D45" "0.00005326"
Set-TextValue "E45" "-0.96%"
Set-TextValue "E46" "0.02%"
Set-TextValue "E47" "1.26%"
Set-TextValue "D48" "0.1353"
Set-TextValue "E48" "-46.29%"
Set-TextValue "E49" "0.02%"
Set-TextValue "E50" "0.02%"
